$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Harmonpreet Singh"
$ws.Range("B2").Value = "Harmonpreet012@gmail.com"
$ws.Range("C2").Value = "popoxoxo"
$ws.Range("D2").Value = "104.jpg"

# Row 3 - C3 ("123") must be stored as text, not a number, so force the
# cell to Text format before assigning the value, then strip the format
# back to Normal so no residual number-format style is left on the cell.
$ws.Range("A3").Value = "Harmonpreet Singh"
$ws.Range("B3").Value = "Harmonpreet012@gmail.com"
$c3 = $ws.Range("C3")
$c3.NumberFormat = "@"
$c3.Value = "123"
$c3.Style = "Normal"
$ws.Range("D3").Value = "53.jpg"

# Row 4
$ws.Range("A4").Value = "papa shango"
$ws.Range("B4").Value = "papa@123"
$ws.Range("C4").Value = "papa"
$ws.Range("D4").Value = "105.jpg"
